# Mark the ".net 3.5 Async" column (F) as "Done" for the methods that were
# just finished - completing rows that previously only had Sync/Async/.Net
# 3.5 marked as Done.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Methods")

$rows = @(38, 43, 44, 81, 82, 83, 84, 88)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 6).Value = "Done"
}

# Scroll the sheet to follow the newly completed rows, matching where the
# author ended up after finishing this pass of edits (top-left cell A76,
# active cell F89).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 76
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F89").Select()
